# Update localization status for the files that have moved from
# "Ready for handoff" into "In Translation" (Report for Archive generation).
# Files 731f11d6-...md and 86eabc96-...md are now "In Translation";
# c8018f2e-...md remains "Ready for handoff".

$wb = $excel.ActiveWorkbook

# --- Overview sheet: File Name / zh-cn / de-de / Latest Handoff Date ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B8").Value = "In Translation"
$overview.Range("C8").Value = "In Translation"
$overview.Range("B9").Value = "In Translation"
$overview.Range("C9").Value = "In Translation"

# --- Per-locale detail sheets: Source File Name / File Extension / Status / ... ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("C8").Value = "In Translation"
    $ws.Range("C9").Value = "In Translation"
}
